{"js": "// Apply the diff: update the date line and all two-digit-by-two-digit\n// multiplication problems in the table to their new values.\nconst replacements = [\n  [\"2024-05-11 Saturday\", \"2024-05-12 Sunday\"],\n  [\"72\u00d750=\", \"27\u00d729=\"],\n  [\"43\u00d757=\", \"24\u00d765=\"],\n  [\"42\u00d746=\", \"42\u00d766=\"],\n  [\"68\u00d712=\", \"73\u00d786=\"],\n  [\"93\u00d736=\", \"70\u00d720=\"],\n  [\"46\u00d789=\", \"73\u00d787=\"],\n  [\"97\u00d759=\", \"11\u00d773=\"],\n  [\"63\u00d719=\", \"40\u00d728=\"],\n  [\"35\u00d768=\", \"77\u00d760=\"],\n  [\"24\u00d753=\", \"43\u00d727=\"],\n  [\"27\u00d773=\", \"75\u00d728=\"],\n  [\"20\u00d797=\", \"71\u00d718=\"],\n  [\"16\u00d760=\", \"92\u00d756=\"],\n  [\"44\u00d760=\", \"72\u00d722=\"],\n  [\"47\u00d793=\", \"12\u00d791=\"],\n  [\"44\u00d748=\", \"75\u00d778=\"],\n  [\"69\u00d792=\", \"82\u00d721=\"],\n  [\"16\u00d722=\", \"81\u00d793=\"],\n  [\"39\u00d721=\", \"14\u00d727=\"],\n  [\"57\u00d750=\", \"96\u00d723=\"],\n  [\"78\u00d712=\", \"86\u00d775=\"],\n  [\"91\u00d777=\", \"22\u00d714=\"],\n  [\"17\u00d730=\", \"75\u00d744=\"],\n  [\"46\u00d741=\", \"72\u00d722=\"],\n  [\"44\u00d795=\", \"47\u00d796=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Apply the diff: update the date line and all two-digit-by-two-digit\n# multiplication problems in the table to their new values.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-05-11 Saturday\", \"2024-05-12 Sunday\"),\n    @(\"72\u00d750=\", \"27\u00d729=\"),\n    @(\"43\u00d757=\", \"24\u00d765=\"),\n    @(\"42\u00d746=\", \"42\u00d766=\"),\n    @(\"68\u00d712=\", \"73\u00d786=\"),\n    @(\"93\u00d736=\", \"70\u00d720=\"),\n    @(\"46\u00d789=\", \"73\u00d787=\"),\n    @(\"97\u00d759=\", \"11\u00d773=\"),\n    @(\"63\u00d719=\", \"40\u00d728=\"),\n    @(\"35\u00d768=\", \"77\u00d760=\"),\n    @(\"24\u00d753=\", \"43\u00d727=\"),\n    @(\"27\u00d773=\", \"75\u00d728=\"),\n    @(\"20\u00d797=\", \"71\u00d718=\"),\n    @(\"16\u00d760=\", \"92\u00d756=\"),\n    @(\"44\u00d760=\", \"72\u00d722=\"),\n    @(\"47\u00d793=\", \"12\u00d791=\"),\n    @(\"44\u00d748=\", \"75\u00d778=\"),\n    @(\"69\u00d792=\", \"82\u00d721=\"),\n    @(\"16\u00d722=\", \"81\u00d793=\"),\n    @(\"39\u00d721=\", \"14\u00d727=\"),\n    @(\"57\u00d750=\", \"96\u00d723=\"),\n    @(\"78\u00d712=\", \"86\u00d775=\"),\n    @(\"91\u00d777=\", \"22\u00d714=\"),\n    @(\"17\u00d730=\", \"75\u00d744=\"),\n    @(\"46\u00d741=\", \"72\u00d722=\"),\n    @(\"44\u00d795=\", \"47\u00d796=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
